$d = $word.ActiveDocument

# 1) "Application Rule" -> "Tax on Subtotal"
$d.Content.Find.Execute("Application Rule", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tax on Subtotal", 2)

# 2) "T" -> "There is certain amount of tax determined according to the subtotal"
$d.Content.Find.Execute("he application must be able to add an item to the cart ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".", 2)

# 3) "and application must be able to check out the items. " -> new text
$d.Content.Find.Execute("and application must be able to check out the items. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " If subtotal goes over 150`$, 20% tax will be added. If it`u{2019}s less than 150 and more than 100, 17% tax will be added to subtotal. If subtotal is more than 50 and less than 100, 13% of tax will be added to subtotal.", 2)

# 4) "T" -> long text (do last since "T" is a very short/ambiguous match)
$d.Content.Find.Execute("Description: T", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Description: There is certain amount of tax determined according to the subtotal", 2)
